$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.004.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.212.80'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.22%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '295.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.508'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.15%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").Value = '  -4.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0772'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.72'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.24'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.91%  '
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.552.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.16'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.230.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.710'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '38.913.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.56%  '
$ws.Range("E20").Value = '  -4.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.38%  '
$ws.Range("E23").Value = '  -5.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '225.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.97%  '
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("E26").Value = '  -7.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.28%  '
$ws.Range("E29").Value = '  -1.45%  '
$ws.Range("E30").Value = '  -2.16%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '147.29'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.79%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.53'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.84%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.82'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.72%  '
$ws.Range("E35").Value = '  -3.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.71%  '
$ws.Range("E37").Value = '  -3.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0949'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.57%  '
$ws.Range("E41").Value = '  -4.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.905.75'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0256'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.22%  '
$ws.Range("E45").Value = '  -16.78%  '
$ws.Range("E46").Value = '  -4.07%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.84%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.425.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.46'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.55%  '
